# Generate Report for Handoff
# Adds a new handoff entry (ab4078aa-f039-478f-8f90-5c6973940d4e) as row 3
# on the Overview, zh-cn and de-de worksheets, extends each table to
# include the new row, and wires up the matching hyperlinks.

$wb = $excel.ActiveWorkbook

$newBase = "ab4078aa-f039-478f-8f90-5c6973940d4e"
$commitSha = "0a357ffb78e8631627c21dd6a4198665ff15ad7f"
$ghBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fbab17eca8384c18d95238340e13f0351d924b59/e2e/"

# BGR encoding of the existing hyperlink font color (RGB FF6495ED)
$hyperlinkColor = 15570276
$dateFormat = "yyyy-mm-dd HH:mm:ss"

function Set-HyperlinkLook($rng) {
    $rng.Font.Color = $hyperlinkColor
    $rng.Font.Underline = 2
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = ($newBase + ".md")
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-29 16:43:48"
$wsOverview.Range("G3").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), ($ghBase + $newBase + ".md"), "", "", ("e2e\" + $newBase + ".md")) | Out-Null
Set-HyperlinkLook $wsOverview.Range("B3")

$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.Resize($wsOverview.Range("A1:G3")) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "'False"
$wsZhCn.Range("G3").Value = ($newBase + "." + $commitSha + ".zh-cn.xlf")
$wsZhCn.Range("H3").Value = "2016-08-29 16:43:43"
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("M3").Value = "'True"
$wsZhCn.Range("O3").Value = "'False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), ($ghBase + $newBase + ".md"), "", "", ($newBase + ".md")) | Out-Null
Set-HyperlinkLook $wsZhCn.Range("A3")

$tblZhCn = $wsZhCn.ListObjects.Item(1)
$tblZhCn.Resize($wsZhCn.Range("A1:P3")) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "'False"
$wsDeDe.Range("G3").Value = ($newBase + "." + $commitSha + ".de-de.xlf")
$wsDeDe.Range("H3").Value = "2016-08-29 16:43:48"
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("M3").Value = "'True"
$wsDeDe.Range("O3").Value = "'False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), ($ghBase + $newBase + ".md"), "", "", ($newBase + ".md")) | Out-Null
Set-HyperlinkLook $wsDeDe.Range("A3")

$tblDeDe = $wsDeDe.ListObjects.Item(1)
$tblDeDe.Resize($wsDeDe.Range("A1:P3")) | Out-Null

Write-Output "Handoff row added for $newBase"
